$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.193.23"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "2.950.74"
$ws.Range("E3").Value = "  -1.61%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'569.28"
$ws.Range("E5").Value = "  -2.75%  "
$ws.Range("D6").Value = "'159.73"
$ws.Range("E6").Value = "  +3.99%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.518"
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("D9").Value = "2.945.03"
$ws.Range("E9").Value = "  -1.60%  "
$ws.Range("D10").Value = "'6.66"
$ws.Range("E10").Value = "  -4.42%  "
$ws.Range("E11").Value = "  -1.31%  "
$ws.Range("E12").Value = "  +1.83%  "
$ws.Range("D13").Value = "'0.0000245"
$ws.Range("E13").Value = "  +1.92%  "
$ws.Range("D14").Value = "'34.12"
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("D16").Value = "65.180.36"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").Value = "3.439.62"
$ws.Range("E17").Value = "  -1.63%  "
$ws.Range("D18").Value = "'6.93"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").Value = "2.946.93"
$ws.Range("E19").Value = "  -1.72%  "
$ws.Range("D20").Value = "'446.64"
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("D21").Value = "'14.48"
$ws.Range("E21").Value = "  +5.75%  "
$ws.Range("D22").Value = "'0.685"
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("E23").Value = "  -1.03%  "
$ws.Range("D24").Value = "'82.27"
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("D25").Value = "'2.21"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("D26").Value = "'12.06"
$ws.Range("E26").Value = "  -3.05%  "
$ws.Range("D27").Value = "'10.06"
$ws.Range("E27").Value = "  -5.85%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "'8.02"
$ws.Range("E29").Value = "  +3.22%  "
$ws.Range("D30").Value = "'2.38"
$ws.Range("E30").Value = "  -0.64%  "
$ws.Range("E31").Value = "  -0.83%  "
$ws.Range("D32").Value = "'0.0000101"
$ws.Range("E32").Value = "  -2.24%  "
$ws.Range("D33").Value = "'27.08"
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("E34").Value = "  -0.73%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("D37").Value = "'5.69"
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("D38").Value = "'49.02"
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("B39").Value = "Arweave"
$ws.Range("C39").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D39").Value = "'44.13"
$ws.Range("E39").Value = "  -3.71%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'1.97"
$ws.Range("E40").Value = "  -6.39%  "
$ws.Range("E41").Value = "  -1.95%  "
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("D43").Value = "'0.299"
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("D44").Value = "'8.42"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").Value = "'385.32"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").Value = "2.714.87"
$ws.Range("E47").Value = "  -1.93%  "
$ws.Range("D48").Value = "'132.99"
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("E50").Value = "  +4.94%  "
$ws.Range("E51").Value = "  +0.58%  "
